$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-14 Thursday", "2024-11-15 Friday"),
    @("53÷3=", "97÷7="),
    @("78÷9=", "86÷3="),
    @("87÷8=", "61÷3="),
    @("64÷8=", "21÷5="),
    @("44÷7=", "43÷8="),
    @("94÷3=", "83÷8="),
    @("18÷2=", "88÷3="),
    @("76÷2=", "58÷7="),
    @("52÷7=", "34÷4="),
    @("79÷9=", "13÷9="),
    @("35÷2=", "21÷8="),
    @("86÷6=", "43÷8="),
    @("79÷6=", "91÷7="),
    @("77÷3=", "63÷5="),
    @("73÷7=", "15÷8="),
    @("92÷7=", "53÷9="),
    @("56÷3=", "29÷5="),
    @("14÷6=", "59÷5="),
    @("47÷4=", "31÷7="),
    @("31÷9=", "95÷5="),
    @("62÷7=", "95÷4="),
    @("62÷8=", "41÷4="),
    @("88÷6=", "33÷4="),
    @("66÷6=", "53÷8="),
    @("48÷3=", "27÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
